$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:K -> E:L.
$ws.Columns("D:D").Insert()

# Carry over number formatting/style from column E (the old column D, now shifted)
# into the freshly inserted (blank) column D.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Populate the new column D with the latest period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 152700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -1300
$ws.Range("D17").Value = 38800
$ws.Range("D18").Value = 113900
$ws.Range("D20").Value = -64400
$ws.Range("D21").Value = 56000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 49500
$ws.Range("D24").Value = 10000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 39500
$ws.Range("D27").Value = 38100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 64400
$ws.Range("D33").Value = 38100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 38100

$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 102800
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 42800
$ws.Range("D49").Value = 76200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 7900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4311700
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 39200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3915400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 17300
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 279900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 379000
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 38100
$ws.Range("D83").Value = 6500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 65100
$ws.Range("D91").Value = -2800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -225400
$ws.Range("D96").Value = -16400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 163800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 3600
